$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'055.348.991-74"
$ws.Range("B2").Value = 3179
$ws.Range("C2").Value = "WP08052004"
$ws.Range("F2").Value = "SPI20250508002150"
$ws.Range("G2").Value = "Teste de registro via automação1"

# Row 3
$ws.Range("A3").Value = "'05534899174"
$ws.Range("B3").Value = 3179
$ws.Range("E3").Value = "Dúvida técnica"
$ws.Range("F3").Value = "'"
$ws.Range("G3").Value = "Teste de registro via automação2"
